# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The error matrix in columns B:K (rows 2-24) is a rolling "age" series:
# column B holds the newest (age-1) error and each column to the right holds
# an older vintage. A new age-1 observation needs to be inserted at the
# front of every row, pushing the existing values one column to the right
# (and, for the fully-populated rows, dropping the oldest value that falls
# out of the 10-column window B:K).

$ws = $excel.ActiveWorkbook.ActiveSheet

# New "age-1" value to insert into column B for each row (2-24).
$newB = @{
    2  = 2.057869132359739
    3  = 6.652313087672924
    4  = -18.36749132628568
    5  = 7.513167073507937
    6  = 0.9564081874156993
    7  = -4.157449276732949
    8  = 1.546611864454844
    9  = 1.156631887942306
    10 = -1.025188112727922
    11 = 0.08364543516793629
    12 = -0.1538585523806955
    13 = 0.7495351060200912
    14 = 0.03849281619118239
    15 = -0.2590580299438133
    16 = 0.01855976243503714
    17 = 0.1467044301255134
    18 = -0.1819613811903656
    19 = 0.4718454808444464
    20 = -0.08594117411414147
    21 = -0.07695400962807622
    22 = -0.5068991247689255
    23 = 0.6215838649243215
    24 = -0.2766911554241067
}

# Columns B..K as 1-based column indices (B=2 ... K=11).
$firstCol = 2   # B
$lastCol  = 11  # K

for ($row = 2; $row -le 24; $row++) {

    # Read the existing row values (B..K) before overwriting anything.
    # NOTE: use .Value2 for reads -- .Value's getter is unreliable here.
    $oldValues = @()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $oldValues += $ws.Cells.Item($row, $col).Value2
    }

    # Shift every existing value one column to the right. The old value that
    # was in column K falls off the end of the 10-column window and is
    # discarded. Walk right-to-left so we never clobber a value before it
    # has been copied.
    for ($col = $lastCol; $col -ge ($firstCol + 1); $col--) {
        $srcIndex = $col - $firstCol - 1
        $ws.Cells.Item($row, $col).Value2 = $oldValues[$srcIndex]
    }

    # Insert the new value at the front (column B).
    $ws.Cells.Item($row, $firstCol).Value2 = $newB[$row]
}
